# Update "想去人数" (desired attendance count) values in column F
# on the "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$exhibitionUpdates = @{
    7  = 11988
    8  = 4454
    9  = 36
    10 = 55
    13 = 2576
    15 = 175
    17 = 5193
    22 = 11420
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new F value
$allTypesUpdates = @{
    7  = 11988
    8  = 4454
    9  = 36
    10 = 55
    13 = 2576
    16 = 175
    18 = 5193
    23 = 11420
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
